$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.203.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.01%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.834.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.91%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +1.21%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'313.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.25%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.97%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4711"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.72%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3693"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.38%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07429"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.65%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.8835"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.19%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'20.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.20%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.838.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.67%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.07332"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.65%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.477"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.98%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'93.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.28%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'6.571"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.91%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'1.013"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.16%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000008822"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.08%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.02%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'14.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.08%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'27.225.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.00%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.308"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.47%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'10.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.07%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.067.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.11%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'1.900"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.15%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'152.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.79%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'18.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.27%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'2.171"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.42%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'5.277"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.08%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'117.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.88%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.08921"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.19%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.7617"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.14%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.176"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.85%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'4.553"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.21%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'2.939"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.06%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.010"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.03%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.103"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.65%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.05337"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.94%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.01959"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.32%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.55%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'7.336"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.15%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'2.399"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.23%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.5342"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.29%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.1666"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.40%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'8.551"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'0.4957"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.10%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'10.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.46%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.010"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.01%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.674"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.08%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'103.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.00%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.06317"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.36%  "
$ws.Range("E51").Style = "Normal"
